$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2535792.38
$ws.Range("C9").Value = 360196.99
$ws.Range("D9").Value = 2895989.37
$ws.Range("E9").Value = 12.43778702129697
$ws.Range("F9").Value = 87.56221297870302
$ws.Range("G9").Value = -65.18875694149187
$ws.Range("H9").Value = -54.20709468988187
$ws.Range("I9").Value = -55.9360167678435
$ws.Range("J9").Value = 25089
$ws.Range("K9").Value = 1042
$ws.Range("L9").Value = 26131
